$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D..M hold per-year (AT_2015..AT_2024) source-column-name mappings.
# Years 2016-2018 (E,F,G) are being brought in line with the other "early/mid"
# years (D=2015, H:K=2019-2022) for several rows, and cleared for field_size.

# Row 3 - farm_id: fill E3:G3 with the same source-column name already used
# in D3/H3:K3 ("hbnr").
$ws.Range("E3:G3").Value = "hbnr"

# Row 4 - crop_code: E4:G4 switch from "SNAR_CODE" to "snart_code" to match
# D4/H4:K4.
$ws.Range("E4:G4").Value = "snart_code"

# Row 5 - crop_name: E5:G5 switch from "SNAR_BEZEICHNUNG" to "snart" to
# match D5/H5:K5.
$ws.Range("E5:G5").Value = "snart"

# Row 11 - organic: fill E11:G11 with "organic" to match D11/H11:K11.
$ws.Range("E11:G11").Value = "organic"

# Row 12 - field_size: clear E12:G12 (was "SL_FLAECHE_BRUTTO_HA").
$ws.Range("E12:G12").ClearContents()

# Update the active selection to match the author's final cursor position.
$ws.Range("F9").Select()
